$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.021462406042487
$block1[0,2] = 1.023215683231386
$block1[0,3] = 1.031418242243197
$block1[0,4] = 1.040297432974917
$block1[1,0] = 1.02
$block1[1,1] = 1.022501836748051
$block1[1,2] = 1.024116155451989
$block1[1,3] = 1.032343811225308
$block1[1,4] = 1.041316064351011
$block1[2,0] = 1.02
$block1[2,1] = 1.023175018296895
$block1[2,2] = 1.02469966425201
$block1[2,3] = 1.032943398160665
$block1[2,4] = 1.041975409810347
$block1[3,0] = 1.02
$block1[3,1] = 1.023458167241409
$block1[3,2] = 1.024945172308576
$block1[3,3] = 1.033195626775223
$block1[3,4] = 1.042252650739759
$block1[4,0] = 1.02
$block1[4,1] = 1.023505717596512
$block1[4,2] = 1.024986405950415
$block1[4,3] = 1.03323798655626
$block1[4,4] = 1.042299203730363
$block1[5,0] = 1.02
$block1[5,1] = 1.023178801183785
$block1[5,2] = 1.024702943954154
$block1[5,3] = 1.032946767815407
$block1[5,4] = 1.041979114114108
$block1[6,0] = 1.02
$block1[6,1] = 1.021813562015933
$block1[6,2] = 1.023519826947792
$block1[6,3] = 1.031730900827825
$block1[6,4] = 1.040641637324521
$block1[7,0] = 1.02
$block1[7,1] = 1.019412447848972
$block1[7,2] = 1.021441518853236
$block1[7,3] = 1.0295936537299
$block1[7,4] = 1.038286593943733
$block1[8,0] = 1.02
$block1[8,1] = 1.017814821620351
$block1[8,2] = 1.020060397508457
$block1[8,3] = 1.028172419600557
$block1[8,4] = 1.036717821552124
$block1[9,0] = 1.02
$block1[9,1] = 1.017123773633112
$block1[9,2] = 1.019463414629235
$block1[9,3] = 1.027557874246489
$block1[9,4] = 1.036038838428846
$block1[10,0] = 1.02
$block1[10,1] = 1.016867198142808
$block1[10,2] = 1.019241827395401
$block1[10,3] = 1.027329734577393
$block1[10,4] = 1.035786680896908
$block1[11,0] = 1.02
$block1[11,1] = 1.016922229454073
$block1[11,2] = 1.019289351411376
$block1[11,3] = 1.027378665421398
$block1[11,4] = 1.035840767397431
$block1[12,0] = 1.02
$block1[12,1] = 1.017102562768148
$block1[12,2] = 1.019445094915672
$block1[12,3] = 1.02753901349657
$block1[12,4] = 1.036017994042766
$block1[13,0] = 1.02
$block1[13,1] = 1.017213686720817
$block1[13,2] = 1.01954107469289
$block1[13,3] = 1.027637826468736
$block1[13,4] = 1.036127195493422
$block1[14,0] = 1.02
$block1[14,1] = 1.017860699688916
$block1[14,2] = 1.020100039521526
$block1[14,3] = 1.028213223111813
$block1[14,4] = 1.036762889960659
$block1[15,0] = 1.02
$block1[15,1] = 1.018266750615706
$block1[15,2] = 1.020450945525342
$block1[15,3] = 1.02857438465939
$block1[15,4] = 1.037161726925901
$block1[16,0] = 1.02
$block1[16,1] = 1.018503664132294
$block1[16,2] = 1.020655724490026
$block1[16,3] = 1.028785126844603
$block1[16,4] = 1.037394391198914
$block1[17,0] = 1.02
$block1[17,1] = 1.018584457538142
$block1[17,2] = 1.020725566005796
$block1[17,3] = 1.028856998495503
$block1[17,4] = 1.037473728659679
$block1[18,0] = 1.02
$block1[18,1] = 1.0182231778616
$block1[18,2] = 1.020413286135118
$block1[18,3] = 1.0285356268762
$block1[18,4] = 1.037118932443042
$block1[19,0] = 1.02
$block1[19,1] = 1.01704945605049
$block1[19,2] = 1.019399227945939
$block1[19,3] = 1.02749179139383
$block1[19,4] = 1.035965803905783
$block1[20,0] = 1.02
$block1[20,1] = 1.016312129950951
$block1[20,2] = 1.018762568395915
$block1[20,3] = 1.026836242230948
$block1[20,4] = 1.03524105888449
$block1[21,0] = 1.02
$block1[21,1] = 1.016702939819327
$block1[21,2] = 1.019099986293642
$block1[21,3] = 1.027183689708057
$block1[21,4] = 1.035625233719083
$block1[22,0] = 1.02
$block1[22,1] = 1.018242866294487
$block1[22,2] = 1.020430302479145
$block1[22,3] = 1.028553139594398
$block1[22,4] = 1.037138269336878
$block1[23,0] = 1.02
$block1[23,1] = 1.020032645285079
$block1[23,2] = 1.021978036406256
$block1[23,3] = 1.030145553445565
$block1[23,4] = 1.038895212904069
$ws.Range("B2:F25").Value2 = $block1

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.026641878174398
$block2[0,1] = 1.02665373659384
$block2[0,2] = 1.026047775447727
$block2[0,3] = 1.034226417757515
$block2[0,4] = 1.043080180560434
$block2[0,5] = 1.012892165284549
$block2[1,0] = 1.026605400056407
$block2[1,1] = 1.027330328448052
$block2[1,2] = 1.026754823851664
$block2[1,3] = 1.03496023264837
$block2[1,4] = 1.043908658534675
$block2[1,5] = 1.013121974758991
$block2[2,0] = 1.026579429592245
$block2[2,1] = 1.027768182983415
$block2[2,2] = 1.027212546304672
$block2[2,3] = 1.035435117346732
$block2[2,4] = 1.04444435664072
$block2[2,5] = 1.013270541973481
$block2[3,0] = 1.026567943967577
$block2[3,1] = 1.027952269425061
$block2[3,2] = 1.027405023583889
$block2[3,3] = 1.035634771790017
$block2[3,4] = 1.044669471310194
$block2[3,5] = 1.013332967025818
$block2[4,0] = 1.026565982163034
$block2[4,1] = 1.027983179057488
$block2[4,2] = 1.027437344322936
$block2[4,3] = 1.035668295382997
$block2[4,4] = 1.044707263575154
$block2[4,5] = 1.013343446545721
$block2[5,0] = 1.026579278352556
$block2[5,1] = 1.027770642707538
$block2[5,2] = 1.0272151179969
$block2[5,3] = 1.035437785089471
$block2[5,4] = 1.044447365000231
$block2[5,5] = 1.013271376228595
$block2[6,0] = 1.026630039631395
$block2[6,1] = 1.026882382349666
$block2[6,2] = 1.02628668087139
$block2[6,3] = 1.034474401555966
$block2[6,4] = 1.043360247038205
$block2[6,5] = 1.012969858212276
$block2[7,0] = 1.026701407511732
$block2[7,1] = 1.025317597656259
$block2[7,2] = 1.02465233131339
$block2[7,3] = 1.032777272441151
$block2[7,4] = 1.041441730556131
$block2[7,5] = 1.012437525677426
$block2[8,0] = 1.026736891602577
$block2[8,1] = 1.024274741406401
$block2[8,2] = 1.023563930099365
$block2[8,3] = 1.031646220315776
$block2[8,4] = 1.040160857393677
$block2[8,5] = 1.012081970157116
$block2[9,0] = 1.026749400744756
$block2[9,1] = 1.023823259265625
$block2[9,2] = 1.023092924476677
$block2[9,3] = 1.031156558926972
$block2[9,4] = 1.039605798009638
$block2[9,5] = 1.011927856130542
$block2[10,0] = 1.026753619156646
$block2[10,1] = 1.023655571451583
$block2[10,2] = 1.022918014550628
$block2[10,3] = 1.030974691369258
$block2[10,4] = 1.03939956031402
$block2[10,5] = 1.011870588138519
$block2[11,0] = 1.02675273364683
$block2[11,1] = 1.023691540471501
$block2[11,2] = 1.02295553138954
$block2[11,3] = 1.031013701910822
$block2[11,4] = 1.039443801874917
$block2[11,5] = 1.011882873361726
$block2[12,0] = 1.026749758166174
$block2[12,1] = 1.023809397875441
$block2[12,2] = 1.023078465496626
$block2[12,3] = 1.031141525398305
$block2[12,4] = 1.039588751635085
$block2[12,5] = 1.011923122811729
$block2[13,0] = 1.026747868184427
$block2[13,1] = 1.023882015417304
$block2[13,2] = 1.023154214910667
$block2[13,3] = 1.031220283608354
$block2[13,4] = 1.039678051528389
$block2[13,5] = 1.011947918762919
$block2[14,0] = 1.026736001330965
$block2[14,1] = 1.024304706549833
$block2[14,2] = 1.023595195124544
$block2[14,3] = 1.031678719521443
$block2[14,4] = 1.040197685831637
$block2[14,5] = 1.012092194931512
$block2[15,0] = 1.026727793721776
$block2[15,1] = 1.024569871709771
$block2[15,2] = 1.02387188554005
$block2[15,3] = 1.031966309299509
$block2[15,4] = 1.040523523595913
$block2[15,5] = 1.012182653953354
$block2[16,0] = 1.026722730762053
$block2[16,1] = 1.024724545890269
$block2[16,2] = 1.024033301286598
$block2[16,3] = 1.032134064345871
$block2[16,4] = 1.040713537413151
$block2[16,5] = 1.0122354020865
$block2[17,0] = 1.026720957646564
$block2[17,1] = 1.02477728710107
$block2[17,2] = 1.024088344431493
$block2[17,3] = 1.032191265951424
$block2[17,4] = 1.040778320084351
$block2[17,5] = 1.012253385269836
$block2[18,0] = 1.026728702819949
$block2[18,1] = 1.024541421184071
$block2[18,2] = 1.02384219648957
$block2[18,3] = 1.031935452723159
$block2[18,4] = 1.040488568606926
$block2[18,5] = 1.012172950114575
$block2[19,0] = 1.026750646179132
$block2[19,1] = 1.02377469145128
$block2[19,2] = 1.023042263287386
$block2[19,3] = 1.031103884161577
$block2[19,4] = 1.039546069299318
$block2[19,5] = 1.011911270989565
$block2[20,0] = 1.026761966814583
$block2[20,1] = 1.023292692088819
$block2[20,2] = 1.022539559905953
$block2[20,3] = 1.030581128157445
$block2[20,4] = 1.038953113262983
$block2[20,5] = 1.011746608958973
$block2[21,0] = 1.026756199900444
$block2[21,1] = 1.023548201880391
$block2[21,2] = 1.02280602893475
$block2[21,3] = 1.030858242790373
$block2[21,4] = 1.039267485115765
$block2[21,5] = 1.011833912023897
$block2[22,0] = 1.026728292889275
$block2[22,1] = 1.024554276727787
$block2[22,2] = 1.023855611610381
$block2[22,3] = 1.031949395454629
$block2[22,4] = 1.040504363391111
$block2[22,5] = 1.012177334908239
$block2[23,0] = 1.02668509339661
$block2[23,1] = 1.025722075805974
$block2[23,2] = 1.025074647324018
$block2[23,3] = 1.033215959728629
$block2[23,4] = 1.041938046479875
$block2[23,5] = 1.012575265316018
$ws.Range("I2:N25").Value2 = $block2
